$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.319.48"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = "'3.688.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.08%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'683.54"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.36%  '
$ws.Range('D6').Value = "'162.52"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.66%  '
$ws.Range('D7').Value = "'3.686.82"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.07%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = "'0.499"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.20%  '
$ws.Range('D10').Value = "'0.148"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.71%  '
$ws.Range('D11').Value = "'7.26"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.99%  '
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('D13').Value = "'0.0000237"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.17%  '
$ws.Range('D14').Value = "'33.61"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.83%  '
$ws.Range('D15').Value = "'4.310.25"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.04%  '
$ws.Range('D16').Value = "'3.690.29"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.48%  '
$ws.Range('D17').Value = "'69.356.48"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').Value = "'16.32"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.13%  '
$ws.Range('D20').Value = "'6.65"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.73%  '
$ws.Range('D21').Value = "'481.82"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').Value = "'9.80"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.89%  '
$ws.Range('D23').Value = "'0.667"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.56%  '
$ws.Range('D24').Value = "'79.96"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.19%  '
$ws.Range('D25').Value = "'3.834.67"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.04%  '
$ws.Range('D26').Value = "'0.0000128"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.90%  '
$ws.Range('D27').Value = "'11.52"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.75%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = "'9.60"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.78%  '
$ws.Range('E30').Value = '  -10.35%  '
$ws.Range('D31').Value = "'2.76"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.51%  '
$ws.Range('D32').Value = "'2.12"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.72%  '
$ws.Range('D33').Value = "'6.83"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.70%  '
$ws.Range('D34').Value = "'27.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.86%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = "'0.164"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.67%  '
$ws.Range('D37').Value = "'3.652.82"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.21%  '
$ws.Range('D38').Value = "'8.58"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.88%  '
$ws.Range('D39').Value = "'6.11"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('D40').Value = "'0.0943"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.13%  '
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('E44').Value = '  -7.56%  '
$ws.Range('D45').Value = "'157.82"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.36%  '
$ws.Range('D46').Value = "'48.17"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').Value = "'2.85"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.70%  '
$ws.Range('E48').Value = '  -13.42%  '
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('D50').Value = "'390.56"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.11%  '
$ws.Range('E51').Value = '  -5.82%  '
